$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 3 (E3, G3, H3)
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Update the active cell / selection to F3
$ws.Range("F3").Select()
